# Insert a new weekly price-report row for "Repollo" (Crespo record, Primera)
# at row 288. All rows from the old 288 onward shift down by one
# (old row 390 -> new row 391), matching the upstream "semanal" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 288..390 down to 289..391, leaving a blank row 288.
$ws.Rows.Item(288).EntireRow.Insert()

# Populate the newly inserted row 288 with this week's record.
$ws.Range("A288").Value = 5
$ws.Range("B288").Value = "Macroferia Regional de Talca"
$ws.Range("C288").Value = "Maule"
$ws.Range("D288").Value = 44837
$ws.Range("E288").Value = 7
$ws.Range("F288").Value = 100112006
$ws.Range("G288").Value = "Repollo"
$ws.Range("H288").Value = "Crespo record"
$ws.Range("I288").Value = "Primera"
$ws.Range("J288").Value = 3000
$ws.Range("K288").Value = 1500
$ws.Range("L288").Value = 1500
$ws.Range("M288").Value = 1500
$ws.Range("N288").Value = "$/unidad"
$ws.Range("O288").Value = "Provincia del Elquí"
$ws.Range("P288").Value = 1500
$ws.Range("Q288").Value = 1
$ws.Range("R288").Value = "Hortaliza"
